$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warning")

# Make "Warning" the active/visible sheet (moves tabSelected off of
# whichever sheet previously had it, e.g. AddOpportunity) and move the
# selection to A11.
$ws.Activate()
$ws.Range("A11").Select()

# Resize the three columns individually (previously B:C shared one <col>
# definition at a wider width).
$ws.Columns.Item(1).ColumnWidth = 34.75
$ws.Columns.Item(2).ColumnWidth = 35.75
$ws.Columns.Item(3).ColumnWidth = 34.75

# Shrink row 2 (it no longer needs to fit the old, much longer message).
$ws.Rows.Item(2).RowHeight = 28.8

# Replace the old "Subject is not an Operating Company" warning copy with
# the new "Companies closed with is missing" warning, and blank out the
# other two warning-message cells on that row.
$ws.Range("A2").Value = "Companies closed with is missing. Please add a counterparty and a closing bid."
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
